# Weekly update to the Haba (Vega Modelo de Temuco) sheet:
# a new price observation is inserted as row 28, pushing the
# existing rows 28-31 down to 29-32.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 28, shifting rows 28:31 down to 29:32.
$ws.Rows("28:28").Insert()

# Populate the newly inserted row 28 with the new observation.
$ws.Range("A28").Value = 10
$ws.Range("B28").Value = "Vega Modelo de Temuco"
$ws.Range("C28").Value = "La Araucanía"
$ws.Range("D28").Value = 44476
$ws.Range("E28").Value = 9
$ws.Range("F28").Value = 100112026
$ws.Range("G28").Value = "Haba"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 50
$ws.Range("K28").Value = 10000
$ws.Range("L28").Value = 10000
$ws.Range("M28").Value = 10000
$ws.Range("N28").Value = "$/saco 25 kilos"
$ws.Range("O28").Value = "Provincia de Limarí"
$ws.Range("P28").Value = 400
$ws.Range("Q28").Value = 25
$ws.Range("R28").Value = "Hortaliza"
